# Update cryptocurrency price/volume figures (and reorder the
# BinanceUSD / Chainlink rows) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.461.77"
$ws.Range("E2").Value = "  +1.80%  "

$ws.Range("D3").Value = "1.827.04"
$ws.Range("E3").Value = "  +1.66%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").Value = "'315.34"
$ws.Range("E5").Value = "  -0.51%  "

$ws.Range("D7").Value = "'0.5059"
$ws.Range("E7").Value = "  -4.90%  "

$ws.Range("D8").Value = "'0.3909"
$ws.Range("E8").Value = "  +0.85%  "

$ws.Range("D9").Value = "'0.07663"
$ws.Range("E9").Value = "  +2.84%  "

$ws.Range("D10").Value = "'41.89"
$ws.Range("E10").Value = "  +1.18%  "

$ws.Range("D11").Value = "'1.108"
$ws.Range("E11").Value = "  +1.85%  "

$ws.Range("D12").Value = "'21.09"
$ws.Range("E12").Value = "  +3.46%  "

$ws.Range("D13").Value = "'6.279"
$ws.Range("E13").Value = "  +1.55%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'7.576"
$ws.Range("E14").Value = "  +1.76%  "

$ws.Range("B15").Value = "BinanceUSD"
$ws.Range("C15").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D15").Value = "'1.001"
$ws.Range("E15").Value = "  +0.25%  "

$ws.Range("D16").Value = "1.823.26"
$ws.Range("E16").Value = "  +1.74%  "

$ws.Range("D17").Value = "'93.10"
$ws.Range("E17").Value = "  +5.25%  "

$ws.Range("E18").Value = "  +2.31%  "

$ws.Range("D19").Value = "'0.06669"
$ws.Range("E19").Value = "  +1.95%  "

$ws.Range("D20").Value = "'17.69"
$ws.Range("E20").Value = "  +2.53%  "

$ws.Range("D22").Value = "'6.149"
$ws.Range("E22").Value = "  +3.12%  "

$ws.Range("D23").Value = "28.495.53"
$ws.Range("E23").Value = "  +1.85%  "

$ws.Range("D24").Value = "'11.12"
$ws.Range("E24").Value = "  +0.07%  "

$ws.Range("D25").Value = "'2.255"
$ws.Range("E25").Value = "  +7.92%  "

$ws.Range("D26").Value = "'156.48"
$ws.Range("E26").Value = "  -0.15%  "

$ws.Range("D27").Value = "'20.60"
$ws.Range("E27").Value = "  +2.25%  "

$ws.Range("D28").Value = "2.036.47"
$ws.Range("E28").Value = "  +1.88%  "

$ws.Range("D29").Value = "'2.397"
$ws.Range("E29").Value = "  +3.87%  "

$ws.Range("D30").Value = "'125.20"
$ws.Range("E30").Value = "  +2.71%  "

$ws.Range("D31").Value = "'1.126"
$ws.Range("E31").Value = "  +2.32%  "

$ws.Range("D32").Value = "'0.1082"
$ws.Range("E32").Value = "  -0.32%  "

$ws.Range("E33").Value = "  +2.98%  "

$ws.Range("D34").Value = "'3.660"

$ws.Range("D35").Value = "'0.07027"
$ws.Range("E35").Value = "  +0.87%  "

$ws.Range("D36").Value = "'0.2224"
$ws.Range("E36").Value = "  +1.00%  "

$ws.Range("D37").Value = "'8.946"
$ws.Range("E37").Value = "  +6.50%  "

$ws.Range("E38").Value = "  +2.15%  "

$ws.Range("D39").Value = "'5.139"
$ws.Range("E39").Value = "  +1.18%  "

$ws.Range("E40").Value = "  +2.08%  "

$ws.Range("D41").Value = "'11.23"
$ws.Range("E41").Value = "  -0.38%  "

$ws.Range("D42").Value = "'1.182"
$ws.Range("E42").Value = "  -0.71%  "

$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = "  +0.16%  "

$ws.Range("D44").Value = "'1.397"

$ws.Range("E45").Value = "  -0.11%  "

$ws.Range("D46").Value = "'0.5901"
$ws.Range("E46").Value = "  +3.22%  "

$ws.Range("E47").Value = "  +1.30%  "

$ws.Range("D48").Value = "'124.42"
$ws.Range("E48").Value = "  -0.08%  "

$ws.Range("D49").Value = "'1.978"
$ws.Range("E49").Value = "  +3.32%  "

$ws.Range("D50").Value = "'1.192"
$ws.Range("E50").Value = "  +1.32%  "

$ws.Range("D51").Value = "'0.06917"
$ws.Range("E51").Value = "  +1.64%  "
